$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8634250164031982
$ws.Range("B1").Value = 1.81453537940979
$ws.Range("D1").Value = 1.914210915565491
$ws.Range("E1").Value = 1.132213473320007
